{"js": "// Update the \"Software\" bullet list entry from \"Godot 4.3.0\" to \"Godot 4.4.1\".\nconst body = context.document.body;\nconst results = body.search(\"Godot 4.3.0\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Godot 4.3.0\" in the document body.');\n}\n\n// Replace the matched text in place, preserving the run's existing formatting.\nresults.items[0].insertText(\"Godot 4.4.1\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the \"Software\" bullet list entry from \"Godot 4.3.0\" to \"Godot 4.4.1\".\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n$found = $range.Find.Execute(\"Godot 4.3.0\")\nif ($found) {\n    $range.Text = \"Godot 4.4.1\"\n}\n"}
